$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Body/Subject/Alias columns (T:V) entirely for all used rows (1-4).
# Clearing the range removes the cells from the sheet data and shrinks the
# sheet dimension back down to A1:S4 automatically.
$ws.Range("T1:V4").Clear()

# Row 2 (S2): shorten the narrative text.
$ws.Range("S2").Value = "Email communication with Imran Rahman regarding Ramiro Gonzalez v. DS Electric, Inc. - Berkley Claim No.: 49538."

# Row 3 (E3/F3): Client/Matter become numeric zero instead of text.
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0

# Row 4 (E4/F4): Client/Matter become numeric zero instead of text.
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

# Row 4 (S4): shorten the narrative text.
$ws.Range("S4").Value = "Email communication with Ms. Yu concerning meeting and conferring on the second and third causes of action in Wang et al. v. Harris et al. (Case No. CIVSB2412923) to avoid a demurrer."
